$d = $word.ActiveDocument

# Paragraph 9
$rng_9_0 = $d.Paragraphs.Item(9).Range
$rng_9_0.Find.Execute("The product vision exists to answer this question.  It helps communicate a sense of direction to stakeholders, both internal and external.  It’s often accompanied with mockups that took a lot of time to create.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_9_0.Text = "A visão do produto existe para responder esta pergunta.  Ela ajuda a comunicar um sentido de direção para stakeholders, tanto internos como externos.  Frequentemente, é acompanhada de mockups que levaram muito tempo para serem criados."

# Paragraph 10
$rng_10_0 = $d.Paragraphs.Item(10).Range
$rng_10_0.Find.Execute("The Issue with Product Visions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_10_0.Text = "O Problema com as Visões de Produto"

# Paragraph 11
$rng_11_0 = $d.Paragraphs.Item(11).Range
$rng_11_0.Find.Execute("But is asking what the product looks like in 3 years even important?  To me, no way.  Product visions are a self-fulfilling prophecy.  They’re created by product leads who will want to make sure they move the product closer to their vision over time so they won’t look like they (a) don’t know how to predict the future, (b) can’t execute or (c) are bad product managers.  So they iterate towards that vision, regardless of whether it’s the right direction.  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_11_0.Text = "Mas perguntar como o produto vai ser em 3 anos é realmente importante?  Para mim, de jeito nenhum.  As visões de produto são uma profecia auto-realizada.  Elas são criadas por líderes de produto que vão querer garantir que o produto seja movido para mais perto da visão ao longo do tempo, para não parecer que (a) não sabem como prever o futuro, (b) não sabem executar ou (c) são gerentes de produto ruins.  Então, iteram em direção a essa visão, independentemente de ser a direção certa ou não. "
$rng_11_1 = $d.Paragraphs.Item(11).Range
$rng_11_1.Find.Execute("At least it’s a direction that stakeholders are familiar with, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_11_1.Text = "Pelo menos é uma direção com que os stakeholders estão familiarizados, "
$rng_11_2 = $d.Paragraphs.Item(11).Range
$rng_11_2.Find.Execute("they think.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_11_2.Text = "pensam."

# Paragraph 12
$rng_12_0 = $d.Paragraphs.Item(12).Range
$rng_12_0.Find.Execute("An Alternative Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_12_0.Text = "Uma Visão Alternativa"

# Paragraph 13
$rng_13_0 = $d.Paragraphs.Item(13).Range
$rng_13_0.Find.Execute("So what is a product lead to communicate if not a product vision? ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_13_0.Text = "Então o que líderes de produto devem comunicar, se não uma visão de produto? "
$rng_13_1 = $d.Paragraphs.Item(13).Range
$rng_13_1.Find.Execute("To me, the best way to communicate a sense of direction to stakeholders is to help them imagine what the customer’s world will be like if the product is successful.  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_13_1.Text = "Para mim, a melhor maneira de comunicar um senso de direção a stakeholders é ajudar a imaginarem como será o mundo dos clientes se o produto for bem sucedido. "
$rng_13_2 = $d.Paragraphs.Item(13).Range
$rng_13_2.Find.Execute("What do I mean by that?  Let’s look at some imaginary examples:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_13_2.Text = "O que quero dizer com isso?  Vamos ver alguns exemplos imaginários:"

# Paragraph 15
$rng_15_0 = $d.Paragraphs.Item(15).Range
$rng_15_0.Find.Execute("Product Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_15_0.Text = "Visão de Produto"
$rng_15_1 = $d.Paragraphs.Item(15).Range
$rng_15_1.Find.Execute("A cross-platform shopping experience that lets customers search, compare and order millions of items by voice.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_15_1.Text = "Uma experiência de compra cross-plataforma que permite aos clientes procurar, comparar e pedir milhões de itens por voz."

# Paragraph 16
$rng_16_0 = $d.Paragraphs.Item(16).Range
$rng_16_0.Find.Execute("Customer Journey Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_16_0.Text = "Visão de Jornada do Cliente"
$rng_16_1 = $d.Paragraphs.Item(16).Range
$rng_16_1.Find.Execute("Imagine never having to leave your house again to go to the store.  No more parking lots, no more lines at the register.  Imagine ordering items from your sofa and opening your front door an hour later to see them there.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_16_1.Text = "Imagine nunca mais precisar sair de casa para ir a uma loja.  Sem mais estacionamentos, sem filas no caixa.  Imagine encomendar produtos do seu sofá, e uma hora depois abrir a porta e encontrá-los lá."

# Paragraph 18
$rng_18_0 = $d.Paragraphs.Item(18).Range
$rng_18_0.Find.Execute("Product Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_18_0.Text = "Visão do produto"
$rng_18_1 = $d.Paragraphs.Item(18).Range
$rng_18_1.Find.Execute("A self-driving car with free WiFi that can be charged in less than 15 minutes.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_18_1.Text = "Um carro autônomo com WiFi gratuito que pode ser carregado em menos de 15 minutos."

# Paragraph 19
$rng_19_0 = $d.Paragraphs.Item(19).Range
$rng_19_0.Find.Execute("Customer Journey Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_19_0.Text = "Visão da Jornada do Cliente"
$rng_19_1 = $d.Paragraphs.Item(19).Range
$rng_19_1.Find.Execute("Imagine being able to check email and read the news while your car drives you to work each morning.  Imagine reducing your carbon footprint and saving money on gas in a luxury, high-tech automobile.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_19_1.Text = "Imagine poder checar o e-mail e ler as notícias enquanto o seu carro te leva para o trabalho todas as manhãs.  Imagine reduzir a sua pegada de carbono e economizar dinheiro com gasolina num carro de alta tecnologia e luxo."

# Paragraph 21
$rng_21_0 = $d.Paragraphs.Item(21).Range
$rng_21_0.Find.Execute("Product Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_21_0.Text = "Visão do Produto"
$rng_21_1 = $d.Paragraphs.Item(21).Range
$rng_21_1.Find.Execute("A flexible subscription, in-home cooking service with a mobile app that uses AI to recommend meals to customers.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_21_1.Text = "Uma serviço de cozinha caseira por assinatura flexível, com um aplicativo que usa IA para recomendar refeições aos clientes."

# Paragraph 22
$rng_22_0 = $d.Paragraphs.Item(22).Range
$rng_22_0.Find.Execute("Customer Journey Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_22_0.Text = "Visão de Jornada do Cliente"
$rng_22_1 = $d.Paragraphs.Item(22).Range
$rng_22_1.Find.Execute("Imagine not having to think about what to make for dinner.  Imagine everything you need to make a gourmet dinner in 30 minutes shows up at your doorstep each week.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_22_1.Text = "Imagine nunca mais precisar pensar no que fazer para o jantar.  Imagine que tudo o que você precisa para fazer um jantar gourmet em 30 minutos aparece na sua porta toda semana."

# Paragraph 23
$rng_23_0 = $d.Paragraphs.Item(23).Range
$rng_23_0.Find.Execute("THE HOMER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_23_0.Text = "O HOMER"

# Paragraph 24
$rng_24_0 = $d.Paragraphs.Item(24).Range
$rng_24_0.Find.Execute("Product Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_24_0.Text = "Visão do Produto"
$rng_24_1 = $d.Paragraphs.Item(24).Range
$rng_24_1.Find.Execute("“Powerful like a gorilla, yet soft and yielding like a Nerf ball.”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_24_1.Text = "“Poderoso como um gorila, mas macio e flexível como uma bola de Nerf.`""
$rng_24_2 = $d.Paragraphs.Item(24).Range
$rng_24_2.Find.Execute("A bubble dome car that can hold huge beverages and plays ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_24_2.Text = "Um carro com teto de bolha que pode segurar bebidas enormes e tocar"
$rng_24_3 = $d.Paragraphs.Item(24).Range
$rng_24_3.Find.Execute(" when you honk.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_24_3.Text = "quando você buzina."

# Paragraph 25
$rng_25_0 = $d.Paragraphs.Item(25).Range
$rng_25_0.Find.Execute("Customer Journey Vision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_25_0.Text = "Visão de Jornada ao Cliente"
$rng_25_1 = $d.Paragraphs.Item(25).Range
$rng_25_1.Find.Execute("Imagine being able to honk many horns when you’re mad.  Imagine being able to shut out screaming kids on road trips with the push of a button.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng_25_1.Text = "Imagine poder usar muitas buzinas quando você está bravo.  Imagine poder isolar crianças berrando em viagens só apertando um botão."
